# Update IG Name / metadata values on the Metadata and Include sheets.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "https://interop.esante.gouv.fr/ig/fhir/tde/ValueSet/MeltingPotVS"
$meta.Range("B3").Value = "2.0.0"
$meta.Range("B8").Value = "2026-01-15T15:25:40+00:00"

# --- Include #1 sheet ---
$inc1 = $wb.Worksheets.Item("Include #1")
$inc1.Range("B4").Value = "https://interop.esante.gouv.fr/ig/fhir/tde/CodeSystem/competence-code-system"

# --- Include #2 sheet ---
$inc2 = $wb.Worksheets.Item("Include #2")
$inc2.Range("B4").Value = "https://interop.esante.gouv.fr/ig/fhir/tde/CodeSystem/type-carte-code-system"
